$wb = $excel.ActiveWorkbook
$wsJournal = $wb.Worksheets.Item("Journal")
$wsTotaux = $wb.Worksheets.Item("Totaux")

# ---------------------------------------------------------------------------
# Journal sheet: add 6 new rows (42-47) describing the Support interview report
# ---------------------------------------------------------------------------
$wsJournal.Rows("42:47").Insert()

# Row 42
$wsJournal.Range("A42").Value = 44992
$wsJournal.Range("B42").Value = 5
$wsJournal.Range("C42").Value = 0.0625
$wsJournal.Range("D42").Value = "Documentation"
$wsJournal.Range("E42").Value = "Rédaction du rapport de projet"

# Row 43
$wsJournal.Range("A43").Value = 44995
$wsJournal.Range("B43").Value = 5
$wsJournal.Range("B43").NumberFormat = "General"
$wsJournal.Range("C43").Value = 0.041666666666666664
$wsJournal.Range("D43").Value = "Entretien"
$wsJournal.Range("E43").Value = "Préparation des documens pour un entretien"

# Row 44
$wsJournal.Range("A44").Value = 44995
$wsJournal.Range("B44").Value = 5
$wsJournal.Range("C44").Value = 0.052083333333333336
$wsJournal.Range("D44").Value = "Entretien"
$wsJournal.Range("E44").Value = "Entretien avec le responsable du support IT"

# Row 45
$wsJournal.Range("A45").Value = 44995
$wsJournal.Range("B45").Value = 5
$wsJournal.Range("C45").Value = 0.020833333333333332
$wsJournal.Range("D45").Value = "Documentation"
$wsJournal.Range("E45").Value = "Mettre les demandes/remarques au propre dans un document Word"

# Row 46
$wsJournal.Range("A46").Value = 44995
$wsJournal.Range("B46").Value = 5
$wsJournal.Range("C46").Value = 0.010416666666666666
$wsJournal.Range("D46").Value = "Documentation"
$wsJournal.Range("E46").Value = "Modification des exports suite à l'entretien "

# Row 47
$wsJournal.Range("A47").Value = 44995
$wsJournal.Range("B47").Value = 5
$wsJournal.Range("C47").Value = 0.03125
$wsJournal.Range("D47").Value = "Documentation"
$wsJournal.Range("E47").Value = "Rédaction du rapport de projet"

# Grow the Journal table (Tableau1) so it covers the new rows
$loJournal = $wsJournal.ListObjects.Item(1)
$loJournal.Resize($wsJournal.Range("A1:F47"))

# ---------------------------------------------------------------------------
# Totaux sheet: add 3 new weekly-total rows (15-17), pushing the Total row to 18
# ---------------------------------------------------------------------------
$wsTotaux.Rows("15:17").Insert()

$wsTotaux.Range("A15").Value = 44991
$wsTotaux.Range("B15").Formula = "=SUM(Journal!C39:C41)"

$wsTotaux.Range("A16").Value = 44992
$wsTotaux.Range("B16").Formula = "=SUM(Journal!C42)"

$wsTotaux.Range("A17").Value = 44995
$wsTotaux.Range("B17").Formula = "=SUM(Journal!C43:C47)"

$wsTotaux.Range("B18").Formula = "=SUM(B2:B17)"

# Grow the Totaux table (Tableau2) so it covers the new rows
$loTotaux = $wsTotaux.ListObjects.Item(1)
$loTotaux.Resize($wsTotaux.Range("A1:B18"))

# ---------------------------------------------------------------------------
# View state: Totaux becomes the active/selected sheet, Journal scrolled down
# ---------------------------------------------------------------------------
$wsJournal.Range("A16").Select()
$wsJournal.Application.ActiveWindow.ScrollRow = 16

$wsTotaux.Activate()
$wsTotaux.Range("I10").Select()

"ok"
